# Append three new match-log rows (5-7) for Mahipal Lomror to the sheet,
# mirroring the existing rows' layout (all values stored as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Rajasthan Royals", "Mumbai Indians", "Mahipal Lomror ", "11", "13", "1", "0", "84.61"),
    @(" Abu Dhabi", " October 03 2020", "RCB won by 8 wickets (with 5 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", "Mahipal Lomror ", "47", "39", "1", "3", "120.51"),
    @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Rajasthan Royals", "Delhi Capitals", "Mahipal Lomror ", "1", "2", "0", "0", "50.00")
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowValues[$c - 1]
    }
}
